$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 13.68414866666667
$ws.Cells.Item(2, 8).Value = 41.052446
$ws.Cells.Item(2, 9).Value = 0.06687192512742757
$ws.Cells.Item(2, 10).Value = 0.06728034099283527
$ws.Cells.Item(2, 13).Value = 7.270285
$ws.Cells.Item(2, 14).Value = 21.810855
$ws.Cells.Item(2, 15).Value = 0.07247008081099036
$ws.Cells.Item(2, 16).Value = 0.07295896149046807
$ws.Cells.Item(2, 17).Value = 99.48766078903667
$ws.Cells.Item(2, 18).Value = 895.3889471013299
$ws.Cells.Item(2, 19).Value = 0.004846213817971173
$ws.Cells.Item(2, 20).Value = 0.004908703807561829

# Row 3
$ws.Cells.Item(3, 7).Value = 13.68414866666667
$ws.Cells.Item(3, 8).Value = 41.052446
$ws.Cells.Item(3, 9).Value = 0.06687192512742757
$ws.Cells.Item(3, 10).Value = 0.06728034099283527
$ws.Cells.Item(3, 15).Value = 0.01763525005644329
$ws.Cells.Item(3, 16).Value = 0.01775421684844728
$ws.Cells.Item(3, 17).Value = 24.20984985681422
$ws.Cells.Item(3, 18).Value = 217.888648711328
$ws.Cells.Item(3, 19).Value = 0.001179303121377939
$ws.Cells.Item(3, 20).Value = 0.001194509763624274

# Row 4
$ws.Cells.Item(4, 7).Value = 13.68414866666667
$ws.Cells.Item(4, 8).Value = 41.052446
$ws.Cells.Item(4, 9).Value = 0.06687192512742757
$ws.Cells.Item(4, 10).Value = 0.06728034099283527
$ws.Cells.Item(4, 13).Value = 31.88752633333333
$ws.Cells.Item(4, 14).Value = 95.66257899999999
$ws.Cells.Item(4, 15).Value = 0.3178543358670601
$ws.Cells.Item(4, 16).Value = 0.3199985703146373
$ws.Cells.Item(4, 17).Value = 436.3536509575815
$ws.Cells.Item(4, 18).Value = 3927.182858618233
$ws.Cells.Item(4, 19).Value = 0.02125553134953026
$ws.Cells.Item(4, 20).Value = 0.02152961292798857

# Row 5
$ws.Cells.Item(5, 7).Value = 13.68414866666667
$ws.Cells.Item(5, 8).Value = 41.052446
$ws.Cells.Item(5, 9).Value = 0.06687192512742757
$ws.Cells.Item(5, 10).Value = 0.06728034099283527
$ws.Cells.Item(5, 13).Value = 2.0166855
$ws.Cells.Item(5, 14).Value = 4.033371
$ws.Cells.Item(5, 15).Value = 0.02010228775836882
$ws.Cells.Item(5, 16).Value = 0.01349193140139489
$ws.Cells.Item(5, 17).Value = 27.596624195911
$ws.Cells.Item(5, 18).Value = 165.579745175466
$ws.Cells.Item(5, 19).Value = 0.001344278681867644
$ws.Cells.Item(5, 20).Value = 0.0009077417453377897

# Row 6
$ws.Cells.Item(6, 7).Value = 13.68414866666667
$ws.Cells.Item(6, 8).Value = 41.052446
$ws.Cells.Item(6, 9).Value = 0.06687192512742757
$ws.Cells.Item(6, 10).Value = 0.06728034099283527
$ws.Cells.Item(6, 13).Value = 57.37750733333333
$ws.Cells.Item(6, 14).Value = 172.132522
$ws.Cells.Item(6, 15).Value = 0.5719380455071374
$ws.Cells.Item(6, 16).Value = 0.5757963199450524
$ws.Cells.Item(6, 17).Value = 785.1623404720901
$ws.Cells.Item(6, 18).Value = 7066.461064248811
$ws.Cells.Item(6, 19).Value = 0.03824659815668056
$ws.Cells.Item(6, 20).Value = 0.0387397727483228

# Row 7
$ws.Cells.Item(7, 9).Value = 0.622425171752627
$ws.Cells.Item(7, 10).Value = 0.6262265923740385
$ws.Cells.Item(7, 13).Value = 7.270285
$ws.Cells.Item(7, 14).Value = 21.810855
$ws.Cells.Item(7, 15).Value = 0.07247008081099036
$ws.Cells.Item(7, 16).Value = 0.07295896149046807
$ws.Cells.Item(7, 17).Value = 926.0033150815517
$ws.Cells.Item(7, 18).Value = 8334.029835733965
$ws.Cells.Item(7, 19).Value = 0.04510720249570744
$ws.Cells.Item(7, 20).Value = 0.04568884183732452

# Row 8
$ws.Cells.Item(8, 9).Value = 0.622425171752627
$ws.Cells.Item(8, 10).Value = 0.6262265923740385
$ws.Cells.Item(8, 15).Value = 0.01763525005644329
$ws.Cells.Item(8, 16).Value = 0.01775421684844728
$ws.Cells.Item(8, 19).Value = 0.01097662354528224
$ws.Cells.Item(8, 20).Value = 0.01111816271727288

# Row 9
$ws.Cells.Item(9, 9).Value = 0.622425171752627
$ws.Cells.Item(9, 10).Value = 0.6262265923740385
$ws.Cells.Item(9, 13).Value = 31.88752633333333
$ws.Cells.Item(9, 14).Value = 95.66257899999999
$ws.Cells.Item(9, 15).Value = 0.3178543358670601
$ws.Cells.Item(9, 16).Value = 0.3199985703146373
$ws.Cells.Item(9, 17).Value = 4061.457713750828
$ws.Cells.Item(9, 18).Value = 36553.11942375746
$ws.Cells.Item(9, 19).Value = 0.1978405395943721
$ws.Cells.Item(9, 20).Value = 0.2003916142526995

# Row 10
$ws.Cells.Item(10, 9).Value = 0.622425171752627
$ws.Cells.Item(10, 10).Value = 0.6262265923740385
$ws.Cells.Item(10, 13).Value = 2.0166855
$ws.Cells.Item(10, 14).Value = 4.033371
$ws.Cells.Item(10, 15).Value = 0.02010228775836882
$ws.Cells.Item(10, 16).Value = 0.01349193140139489
$ws.Cells.Item(10, 17).Value = 256.8616578960655
$ws.Cells.Item(10, 18).Value = 1541.169947376393
$ws.Cells.Item(10, 19).Value = 0.01251216991062345
$ws.Cells.Item(10, 20).Value = 0.008449006226039807

# Row 11
$ws.Cells.Item(11, 9).Value = 0.622425171752627
$ws.Cells.Item(11, 10).Value = 0.6262265923740385
$ws.Cells.Item(11, 13).Value = 57.37750733333333
$ws.Cells.Item(11, 14).Value = 172.132522
$ws.Cells.Item(11, 15).Value = 0.5719380455071374
$ws.Cells.Item(11, 16).Value = 0.5757963199450524
$ws.Cells.Item(11, 17).Value = 7308.071416977837
$ws.Cells.Item(11, 18).Value = 65772.64275280053
$ws.Cells.Item(11, 19).Value = 0.3559886362066418
$ws.Cells.Item(11, 20).Value = 0.3605789673407018

# Row 12
$ws.Cells.Item(12, 7).Value = 17.548286
$ws.Cells.Item(12, 8).Value = 52.644858
$ws.Cells.Item(12, 9).Value = 0.08575525566783661
$ws.Cells.Item(12, 10).Value = 0.08627900022715801
$ws.Cells.Item(12, 13).Value = 7.270285
$ws.Cells.Item(12, 14).Value = 21.810855
$ws.Cells.Item(12, 15).Value = 0.07247008081099036
$ws.Cells.Item(12, 16).Value = 0.07295896149046807
$ws.Cells.Item(12, 17).Value = 127.58104048151
$ws.Cells.Item(12, 18).Value = 1148.22936433359
$ws.Cells.Item(12, 19).Value = 0.006214690308215258
$ws.Cells.Item(12, 20).Value = 0.006294826255009307

# Row 13
$ws.Cells.Item(13, 7).Value = 17.548286
$ws.Cells.Item(13, 8).Value = 52.644858
$ws.Cells.Item(13, 9).Value = 0.08575525566783661
$ws.Cells.Item(13, 10).Value = 0.08627900022715801
$ws.Cells.Item(13, 15).Value = 0.01763525005644329
$ws.Cells.Item(13, 16).Value = 0.01775421684844728
$ws.Cells.Item(13, 17).Value = 31.04624040948267
$ws.Cells.Item(13, 18).Value = 279.416163685344
$ws.Cells.Item(13, 19).Value = 0.001512315377356524
$ws.Cells.Item(13, 20).Value = 0.001531816079500196

# Row 14
$ws.Cells.Item(14, 7).Value = 17.548286
$ws.Cells.Item(14, 8).Value = 52.644858
$ws.Cells.Item(14, 9).Value = 0.08575525566783661
$ws.Cells.Item(14, 10).Value = 0.08627900022715801
$ws.Cells.Item(14, 13).Value = 31.88752633333333
$ws.Cells.Item(14, 14).Value = 95.66257899999999
$ws.Cells.Item(14, 15).Value = 0.3178543358670601
$ws.Cells.Item(14, 16).Value = 0.3199985703146373
$ws.Cells.Item(14, 17).Value = 559.5714319298646
$ws.Cells.Item(14, 18).Value = 5036.142887368782
$ws.Cells.Item(14, 19).Value = 0.02725767983741015
$ws.Cells.Item(14, 20).Value = 0.02760915672086683

# Row 15
$ws.Cells.Item(15, 7).Value = 17.548286
$ws.Cells.Item(15, 8).Value = 52.644858
$ws.Cells.Item(15, 9).Value = 0.08575525566783661
$ws.Cells.Item(15, 10).Value = 0.08627900022715801
$ws.Cells.Item(15, 13).Value = 2.0166855
$ws.Cells.Item(15, 14).Value = 4.033371
$ws.Cells.Item(15, 15).Value = 0.02010228775836882
$ws.Cells.Item(15, 16).Value = 0.01349193140139489
$ws.Cells.Item(15, 17).Value = 35.389373926053
$ws.Cells.Item(15, 18).Value = 212.336243556318
$ws.Cells.Item(15, 19).Value = 0.001723876826227341
$ws.Cells.Item(15, 20).Value = 0.00116407035244575

# Row 16
$ws.Cells.Item(16, 7).Value = 17.548286
$ws.Cells.Item(16, 8).Value = 52.644858
$ws.Cells.Item(16, 9).Value = 0.08575525566783661
$ws.Cells.Item(16, 10).Value = 0.08627900022715801
$ws.Cells.Item(16, 13).Value = 57.37750733333333
$ws.Cells.Item(16, 14).Value = 172.132522
$ws.Cells.Item(16, 15).Value = 0.5719380455071374
$ws.Cells.Item(16, 16).Value = 0.5757963199450524
$ws.Cells.Item(16, 17).Value = 1006.876908652431
$ws.Cells.Item(16, 18).Value = 9061.892177871876
$ws.Cells.Item(16, 19).Value = 0.04904669331862734
$ws.Cells.Item(16, 20).Value = 0.04967913081933593

# Row 17
$ws.Cells.Item(17, 7).Value = 3.726573
$ws.Cells.Item(17, 8).Value = 7.453145999999999
$ws.Cells.Item(17, 9).Value = 0.01821107887003078
$ws.Cells.Item(17, 10).Value = 0.01221486788751604
$ws.Cells.Item(17, 13).Value = 7.270285
$ws.Cells.Item(17, 14).Value = 21.810855
$ws.Cells.Item(17, 15).Value = 0.07247008081099036
$ws.Cells.Item(17, 16).Value = 0.07295896149046807
$ws.Cells.Item(17, 17).Value = 27.093247783305
$ws.Cells.Item(17, 18).Value = 162.55948669983
$ws.Cells.Item(17, 19).Value = 0.001319758357366449
$ws.Cells.Item(17, 20).Value = 0.0008911840758164376

# Row 18
$ws.Cells.Item(18, 7).Value = 3.726573
$ws.Cells.Item(18, 8).Value = 7.453145999999999
$ws.Cells.Item(18, 9).Value = 0.01821107887003078
$ws.Cells.Item(18, 10).Value = 0.01221486788751604
$ws.Cells.Item(18, 15).Value = 0.01763525005644329
$ws.Cells.Item(18, 16).Value = 0.01775421684844728
$ws.Cells.Item(18, 17).Value = 6.593013201488
$ws.Cells.Item(18, 18).Value = 39.55807920892799
$ws.Cells.Item(18, 19).Value = 0.0003211569296706035
$ws.Cells.Item(18, 20).Value = 0.0002168654132500949

# Row 19
$ws.Cells.Item(19, 7).Value = 3.726573
$ws.Cells.Item(19, 8).Value = 7.453145999999999
$ws.Cells.Item(19, 9).Value = 0.01821107887003078
$ws.Cells.Item(19, 10).Value = 0.01221486788751604
$ws.Cells.Item(19, 13).Value = 31.88752633333333
$ws.Cells.Item(19, 14).Value = 95.66257899999999
$ws.Cells.Item(19, 15).Value = 0.3178543358670601
$ws.Cells.Item(19, 16).Value = 0.3199985703146373
$ws.Cells.Item(19, 17).Value = 118.831194670589
$ws.Cells.Item(19, 18).Value = 712.9871680235339
$ws.Cells.Item(19, 19).Value = 0.005788470379656284
$ws.Cells.Item(19, 20).Value = 0.003908740260587306

# Row 20
$ws.Cells.Item(20, 7).Value = 3.726573
$ws.Cells.Item(20, 8).Value = 7.453145999999999
$ws.Cells.Item(20, 9).Value = 0.01821107887003078
$ws.Cells.Item(20, 10).Value = 0.01221486788751604
$ws.Cells.Item(20, 13).Value = 2.0166855
$ws.Cells.Item(20, 14).Value = 4.033371
$ws.Cells.Item(20, 15).Value = 0.02010228775836882
$ws.Cells.Item(20, 16).Value = 0.01349193140139489
$ws.Cells.Item(20, 17).Value = 7.515325733791499
$ws.Cells.Item(20, 18).Value = 30.061302935166
$ws.Cells.Item(20, 19).Value = 0.0003660843478357088
$ws.Cells.Item(20, 20).Value = 0.0001648021596154677

# Row 21
$ws.Cells.Item(21, 7).Value = 3.726573
$ws.Cells.Item(21, 8).Value = 7.453145999999999
$ws.Cells.Item(21, 9).Value = 0.01821107887003078
$ws.Cells.Item(21, 10).Value = 0.01221486788751604
$ws.Cells.Item(21, 13).Value = 57.37750733333333
$ws.Cells.Item(21, 14).Value = 172.132522
$ws.Cells.Item(21, 15).Value = 0.5719380455071374
$ws.Cells.Item(21, 16).Value = 0.5757963199450524
$ws.Cells.Item(21, 17).Value = 213.821469635702
$ws.Cells.Item(21, 18).Value = 1282.928817814212
$ws.Cells.Item(21, 19).Value = 0.01041560885550173
$ws.Cells.Item(21, 20).Value = 0.00703327597824673

# Row 22
$ws.Cells.Item(22, 7).Value = 42.30495733333333
$ws.Cells.Item(22, 8).Value = 126.914872
$ws.Cells.Item(22, 9).Value = 0.206736568582078
$ws.Cells.Item(22, 10).Value = 0.2079991985184523
$ws.Cells.Item(22, 13).Value = 7.270285
$ws.Cells.Item(22, 14).Value = 21.810855
$ws.Cells.Item(22, 15).Value = 0.07247008081099036
$ws.Cells.Item(22, 16).Value = 0.07295896149046807
$ws.Cells.Item(22, 17).Value = 307.5690967261734
$ws.Cells.Item(22, 18).Value = 2768.12187053556
$ws.Cells.Item(22, 19).Value = 0.01498221583173004
$ws.Cells.Item(22, 20).Value = 0.01517540551475598

# Row 23
$ws.Cells.Item(23, 7).Value = 42.30495733333333
$ws.Cells.Item(23, 8).Value = 126.914872
$ws.Cells.Item(23, 9).Value = 0.206736568582078
$ws.Cells.Item(23, 10).Value = 0.2079991985184523
$ws.Cells.Item(23, 15).Value = 0.01763525005644329
$ws.Cells.Item(23, 16).Value = 0.01775421684844728
$ws.Cells.Item(23, 17).Value = 74.84547926125511
$ws.Cells.Item(23, 18).Value = 673.609313351296
$ws.Cells.Item(23, 19).Value = 0.003645851082755984
$ws.Cells.Item(23, 20).Value = 0.003692862874799836

# Row 24
$ws.Cells.Item(24, 7).Value = 42.30495733333333
$ws.Cells.Item(24, 8).Value = 126.914872
$ws.Cells.Item(24, 9).Value = 0.206736568582078
$ws.Cells.Item(24, 10).Value = 0.2079991985184523
$ws.Cells.Item(24, 13).Value = 31.88752633333333
$ws.Cells.Item(24, 14).Value = 95.66257899999999
$ws.Cells.Item(24, 15).Value = 0.3178543358670601
$ws.Cells.Item(24, 16).Value = 0.3199985703146373
$ws.Cells.Item(24, 17).Value = 1349.00044099721
$ws.Cells.Item(24, 18).Value = 12141.00396897489
$ws.Cells.Item(24, 19).Value = 0.06571211470609134
$ws.Cells.Item(24, 20).Value = 0.06655944615249515

# Row 25
$ws.Cells.Item(25, 7).Value = 42.30495733333333
$ws.Cells.Item(25, 8).Value = 126.914872
$ws.Cells.Item(25, 9).Value = 0.206736568582078
$ws.Cells.Item(25, 10).Value = 0.2079991985184523
$ws.Cells.Item(25, 13).Value = 2.0166855
$ws.Cells.Item(25, 14).Value = 4.033371
$ws.Cells.Item(25, 15).Value = 0.02010228775836882
$ws.Cells.Item(25, 16).Value = 0.01349193140139489
$ws.Cells.Item(25, 17).Value = 85.315794032252
$ws.Cells.Item(25, 18).Value = 511.894764193512
$ws.Cells.Item(25, 19).Value = 0.004155877991814683
$ws.Cells.Item(25, 20).Value = 0.002806310917956075

# Row 26
$ws.Cells.Item(26, 7).Value = 42.30495733333333
$ws.Cells.Item(26, 8).Value = 126.914872
$ws.Cells.Item(26, 9).Value = 0.206736568582078
$ws.Cells.Item(26, 10).Value = 0.2079991985184523
$ws.Cells.Item(26, 13).Value = 57.37750733333333
$ws.Cells.Item(26, 14).Value = 172.132522
$ws.Cells.Item(26, 15).Value = 0.5719380455071374
$ws.Cells.Item(26, 16).Value = 0.5757963199450524
$ws.Cells.Item(26, 17).Value = 2427.352999629687
$ws.Cells.Item(26, 18).Value = 21846.17699666719
$ws.Cells.Item(26, 19).Value = 0.118240508969686
$ws.Cells.Item(26, 20).Value = 0.1197651730584452
